# Change description:
#  - Sheet1: replace the repeated "hello" text in A1:A11 with distinct
#    strings (aaa, bbb, ccc, ... kkk) while keeping existing cell styles,
#    and move the active selection to C6.
#  - Add a new worksheet "Sheet2" after Sheet1 containing a small table
#    (A1:C2) with a mix of plain and formatted (italic / bold) cells,
#    and make it the active sheet with C1 selected.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: give every row its own distinct text -------------------------
$ws1.Range("A1").Value  = "aaa"
$ws1.Range("A2").Value  = "bbb"
$ws1.Range("A3").Value  = "ccc"
$ws1.Range("A4").Value  = "ddd"
$ws1.Range("A5").Value  = "eee"
$ws1.Range("A6").Value  = "fff"
$ws1.Range("A7").Value  = "ggg"
$ws1.Range("A8").Value  = "hhh"
$ws1.Range("A9").Value  = "iii"
$ws1.Range("A10").Value = "jjj"
$ws1.Range("A11").Value = "kkk"

# Keep Sheet1's selection parked on C6 (it will lose tabSelected once
# Sheet2 below becomes the active sheet).
$ws1.Range("C6").Select() | Out-Null

# --- Add Sheet2 right after Sheet1 -----------------------------------------
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "llll"
$ws2.Range("B1").Value = "mmmm"
$ws2.Range("C1").Value = "nnnn"
$ws2.Range("A2").Value = "oooo"
$ws2.Range("B2").Value = "pppp"

# Formatting that distinguishes individual cells in the dataframe
$ws2.Range("C1").Font.Italic = $true
$ws2.Range("B2").Font.Bold   = $true

# Sheet2 ends up the active sheet with C1 selected
$ws2.Range("C1").Select() | Out-Null
